$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1005.05554
$ws.Range("I18").Value = 711.2353000000001
$ws.Range("J18").Value = 6000
$ws.Range("K18").Value = 711.2353000000001
$ws.Range("L18").Value = 6000
$ws.Range("M18").Value = -427.2353000000001
$ws.Range("N18").Value = -6568

$ws.Range("H51").Value = 7680.857
$ws.Range("I51").Value = 14750
$ws.Range("J51").Value = 3330.6155
$ws.Range("K51").Value = 14750
$ws.Range("L51").Value = 3330.6155
$ws.Range("M51").Value = -14266
$ws.Range("N51").Value = -4298.6155

$ws.Range("H107").Value = 461.25
$ws.Range("I107").Value = 447.4737
$ws.Range("J107").Value = 513.6
$ws.Range("K107").Value = 447.4737
$ws.Range("L107").Value = 513.6
$ws.Range("M107").Value = 1472.5263
$ws.Range("N107").Value = -4353.6

$ws.Range("H112").Value = 1423
$ws.Range("I112").Value = 490
$ws.Range("J112").Value = 1609.6
$ws.Range("K112").Value = 1470
$ws.Range("L112").Value = 4828.799999999999
$ws.Range("M112").Value = -362
$ws.Range("N112").Value = -7044.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27782.412
$ws.Range("I32").Value = 5752.9688
$ws.Range("J32").Value = 380253.5
$ws.Range("K32").Value = 5752.9688
$ws.Range("L32").Value = 380253.5
$ws.Range("M32").Value = -5465.9688
$ws.Range("N32").Value = -380827.5

$ws.Range("H45").Value = 112688.336
$ws.Range("I45").Value = 144314
$ws.Range("J45").Value = 1998.5
$ws.Range("K45").Value = 144314
$ws.Range("L45").Value = 1998.5
$ws.Range("M45").Value = -143937
$ws.Range("N45").Value = -2752.5

$ws.Range("H82").Value = 30567
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 30567
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 30567
$ws.Range("N82").Value = -31289

$ws.Range("H85").Value = 30567
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 30567
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 30567
$ws.Range("N85").Value = -33063

$ws.Range("H122").Value = 1743.6086
$ws.Range("I122").Value = 1710.3125
$ws.Range("J122").Value = 1819.7142
$ws.Range("K122").Value = 5130.9375
$ws.Range("L122").Value = 5459.142599999999
$ws.Range("M122").Value = -2680.9375
$ws.Range("N122").Value = -10359.1426

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 102302.73
$ws.Range("I86").Value = 159654.86
$ws.Range("J86").Value = 1936.5
$ws.Range("K86").Value = 159654.86
$ws.Range("L86").Value = 1936.5
$ws.Range("M86").Value = -158531.86
$ws.Range("N86").Value = -4182.5

$ws.Range("H89").Value = 102302.73
$ws.Range("I89").Value = 159654.86
$ws.Range("J89").Value = 1936.5
$ws.Range("K89").Value = 798274.2999999999
$ws.Range("L89").Value = 9682.5
$ws.Range("M89").Value = -792658.2999999999
$ws.Range("N89").Value = -20914.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 690
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 690
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 690
$ws.Range("N13").Value = -968

$ws.Range("H31").Value = 19935.424
$ws.Range("I31").Value = 1315.8334
$ws.Range("J31").Value = 35895.07
$ws.Range("K31").Value = 1315.8334
$ws.Range("L31").Value = 35895.07
$ws.Range("M31").Value = -1020.8334
$ws.Range("N31").Value = -36485.07

$ws.Range("H34").Value = 19935.424
$ws.Range("I34").Value = 1315.8334
$ws.Range("J34").Value = 35895.07
$ws.Range("K34").Value = 1315.8334
$ws.Range("L34").Value = 35895.07
$ws.Range("M34").Value = -1113.8334
$ws.Range("N34").Value = -36299.07

$ws.Range("H41").Value = 9330.143
$ws.Range("I41").Value = 4183.6665
$ws.Range("J41").Value = 13190
$ws.Range("K41").Value = 4183.6665
$ws.Range("L41").Value = 13190
$ws.Range("M41").Value = -3755.6665
$ws.Range("N41").Value = -14046

$ws.Range("H51").Value = 7955.3335
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 7955.3335
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 7955.3335
$ws.Range("N51").Value = -9427.333500000001

$ws.Range("H59").Value = 28990
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 28990
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 28990
$ws.Range("N59").Value = -31280

$ws.Range("H60").Value = 12820.571
$ws.Range("I60").Value = 4868
$ws.Range("J60").Value = 18785
$ws.Range("K60").Value = 4868
$ws.Range("L60").Value = 18785
$ws.Range("M60").Value = -4357
$ws.Range("N60").Value = -19807

$ws.Range("H61").Value = 7955.3335
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 7955.3335
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 7955.3335
$ws.Range("N61").Value = -8651.333500000001

$ws.Range("H68").Value = 18367.264
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 18367.264
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 18367.264
$ws.Range("N68").Value = -19865.264

$ws.Range("H71").Value = 18367.264
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 18367.264
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 55101.792
$ws.Range("N71").Value = -62589.792

$ws.Range("H99").Value = 10604.385
$ws.Range("I99").Value = 3559.25
$ws.Range("J99").Value = 21876.6
$ws.Range("K99").Value = 3559.25
$ws.Range("L99").Value = 21876.6
$ws.Range("M99").Value = -2061.25
$ws.Range("N99").Value = -24872.6

$ws.Range("H105").Value = 1214.0588
$ws.Range("I105").Value = 1178.091
$ws.Range("J105").Value = 1280
$ws.Range("K105").Value = 1178.091
$ws.Range("L105").Value = 1280
$ws.Range("M105").Value = 568.9090000000001
$ws.Range("N105").Value = -4774

$ws.Range("H126").Value = 10604.385
$ws.Range("I126").Value = 3559.25
$ws.Range("J126").Value = 21876.6
$ws.Range("K126").Value = 10677.75
$ws.Range("L126").Value = 65629.79999999999
$ws.Range("M126").Value = -8207.75
$ws.Range("N126").Value = -70569.79999999999

$ws.Range("H132").Value = 5600.727
$ws.Range("I132").Value = 6167
$ws.Range("J132").Value = 4921.2
$ws.Range("K132").Value = 18501
$ws.Range("L132").Value = 14763.6
$ws.Range("M132").Value = -15971
$ws.Range("N132").Value = -19823.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 193.81818
$ws.Range("I6").Value = 119
$ws.Range("J6").Value = 393.33334
$ws.Range("K6").Value = 357
$ws.Range("L6").Value = 1180.00002
$ws.Range("M6").Value = -244
$ws.Range("N6").Value = -1406.00002

$ws.Range("H15").Value = 92.95238000000001
$ws.Range("I15").Value = 30
$ws.Range("J15").Value = 294.4
$ws.Range("K15").Value = 90
$ws.Range("L15").Value = 883.1999999999999
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = -1163.2

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("M123").ClearContents()

$ws.Range("H126").Value = 1975.4546
$ws.Range("I126").Value = 1707.5
$ws.Range("J126").Value = 2128.5715
$ws.Range("K126").Value = 5122.5
$ws.Range("L126").Value = 6385.7145
$ws.Range("M126").Value = -182.5
$ws.Range("N126").Value = -16265.7145

$ws.Range("H129").Value = 560437.75
$ws.Range("I129").Value = 11971.2
$ws.Range("J129").Value = 1246021
$ws.Range("K129").Value = 35913.60000000001
$ws.Range("L129").Value = 3738063
$ws.Range("M129").Value = -30913.60000000001
$ws.Range("N129").Value = -3748063

$ws.Range("H130").Value = 330
$ws.Range("I130").Value = 330
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 990
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = 4030

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 167.66667
$ws.Range("I3").Value = 167.66667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 167.66667
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -51.66667000000001

$ws.Range("H102").Value = 252417.5
$ws.Range("I102").Value = 1659.3846
$ws.Range("J102").Value = 548768
$ws.Range("K102").Value = 1659.3846
$ws.Range("L102").Value = 548768
$ws.Range("M102").Value = -37.38460000000009
$ws.Range("N102").Value = -552012

$ws.Range("H113").Value = 1598.1666
$ws.Range("I113").Value = 1650
$ws.Range("J113").Value = 1587.8
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 1587.8
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -5927.8

$ws.Range("H122").Value = 1471.8
$ws.Range("I122").Value = 980
$ws.Range("J122").Value = 1799.6666
$ws.Range("K122").Value = 2940
$ws.Range("L122").Value = 5398.9998
$ws.Range("M122").Value = -490
$ws.Range("N122").Value = -10298.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4270.8
$ws.Range("I7").Value = 2900.8
$ws.Range("J7").Value = 5640.8
$ws.Range("K7").Value = 2900.8
$ws.Range("L7").Value = 5640.8
$ws.Range("M7").Value = -2788.8
$ws.Range("N7").Value = -5864.8

$ws.Range("H40").Value = 113333.336
$ws.Range("I40").Value = 1000000
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 1000000
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -999864
$ws.Range("N40").Value = -2772

$ws.Range("H99").Value = 31900
$ws.Range("I99").Value = 29800
$ws.Range("J99").Value = 34000
$ws.Range("K99").Value = 29800
$ws.Range("L99").Value = 34000
$ws.Range("M99").Value = -26805
$ws.Range("N99").Value = -39990

$ws.Range("H100").Value = 2780
$ws.Range("I100").Value = 2300
$ws.Range("J100").Value = 3740
$ws.Range("K100").Value = 2300
$ws.Range("L100").Value = 3740
$ws.Range("M100").Value = -1759
$ws.Range("N100").Value = -4822

$ws.Range("H122").Value = 2476.5789
$ws.Range("I122").Value = 2430.7334
$ws.Range("J122").Value = 2648.5
$ws.Range("K122").Value = 7292.2002
$ws.Range("L122").Value = 7945.5
$ws.Range("M122").Value = -4842.2002
$ws.Range("N122").Value = -12845.5

$ws.Range("H126").Value = 4270.8
$ws.Range("I126").Value = 2900.8
$ws.Range("J126").Value = 5640.8
$ws.Range("K126").Value = 8702.400000000001
$ws.Range("L126").Value = 16922.4
$ws.Range("M126").Value = -6232.400000000001
$ws.Range("N126").Value = -21862.4

$ws.Range("H127").Value = 29800
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 29800
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 29800
$ws.Range("N127").Value = -39720

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1516
$ws.Range("I122").Value = 871.6667
$ws.Range("J122").Value = 3019.4443
$ws.Range("K122").Value = 2615.0001
$ws.Range("L122").Value = 9058.332900000001
$ws.Range("M122").Value = -165.0001000000002
$ws.Range("N122").Value = -13958.3329
